$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.268596529960632
$ws.Range("B1").Value = 2.743371725082397
$ws.Range("C1").Value = 4.985915184020996
$ws.Range("D1").Value = 2.020536422729492
$ws.Range("E1").Value = 1.031789779663086
